$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H6").Value = 2562.75
$ws_ALC.Range("I6").Value = 2562.75
$ws_ALC.Range("K6").Value = 7688.25
$ws_ALC.Range("M6").Value = -7576.25
$ws_ALC.Range("H62").Value = 7158
$ws_ALC.Range("I62").Value = 6408.727
$ws_ALC.Range("K62").Value = 6408.727
$ws_ALC.Range("M62").Value = -5784.727
$ws_ALC.Range("H65").Value = 7158
$ws_ALC.Range("I65").Value = 6408.727
$ws_ALC.Range("K65").Value = 32043.635
$ws_ALC.Range("M65").Value = -28923.635
$ws_ALC.Range("H98").Value = 2986.85
$ws_ALC.Range("I98").Value = 3072.4119
$ws_ALC.Range("K98").Value = 3072.4119
$ws_ALC.Range("M98").Value = -1574.4119
$ws_ALC.Range("H99").Value = 1046.05
$ws_ALC.Range("I99").Value = 949.8570999999999
$ws_ALC.Range("J99").Value = 1270.5
$ws_ALC.Range("K99").Value = 2849.5713
$ws_ALC.Range("L99").Value = 3811.5
$ws_ALC.Range("M99").Value = -1351.5713
$ws_ALC.Range("N99").Value = -6807.5
$ws_ALC.Range("H100").Value = 1724.75
$ws_ALC.Range("I100").Value = 1724.75
$ws_ALC.Range("K100").Value = 1724.75
$ws_ALC.Range("M100").Value = -1183.75
$ws_ALC.Range("H106").Value = 4449882.5
$ws_ALC.Range("I106").Value = 6066203.5
$ws_ALC.Range("K106").Value = 6066203.5
$ws_ALC.Range("M106").Value = -6065572.5
$ws_ALC.Range("H122").Value = 2986.85
$ws_ALC.Range("I122").Value = 3072.4119
$ws_ALC.Range("K122").Value = 9217.235700000001
$ws_ALC.Range("M122").Value = -6767.235700000001
$ws_ALC.Range("H138").Value = 3228.122
$ws_ALC.Range("I138").Value = 2357.9
$ws_ALC.Range("J138").Value = 4056.9048
$ws_ALC.Range("K138").Value = 7073.700000000001
$ws_ALC.Range("L138").Value = 12170.7144
$ws_ALC.Range("M138").Value = -1933.700000000001
$ws_ALC.Range("N138").Value = -22450.7144

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H4").Value = 33197.6
$ws_ARM.Range("I4").Value = 33197.6
$ws_ARM.Range("K4").Value = 33197.6
$ws_ARM.Range("M4").Value = -33081.6
$ws_ARM.Range("H32").Value = 8241498
$ws_ARM.Range("J32").Value = 71442570
$ws_ARM.Range("L32").Value = 71442570
$ws_ARM.Range("N32").Value = -71443144
$ws_ARM.Range("H61").Value = 3649.7778
$ws_ARM.Range("I61").Value = 3430.7144
$ws_ARM.Range("K61").Value = 3430.7144
$ws_ARM.Range("M61").Value = -3218.7144
$ws_ARM.Range("H70").Value = 85000
$ws_ARM.Range("J70").Value = 85000
$ws_ARM.Range("L70").Value = 85000
$ws_ARM.Range("N70").Value = -85540
$ws_ARM.Range("H73").Value = 85000
$ws_ARM.Range("J73").Value = 85000
$ws_ARM.Range("L73").Value = 85000
$ws_ARM.Range("N73").Value = -86872
$ws_ARM.Range("H122").Value = 3892.745
$ws_ARM.Range("I122").Value = 2819.7576
$ws_ARM.Range("K122").Value = 8459.272799999999
$ws_ARM.Range("M122").Value = -6009.272799999999
$ws_ARM.Range("H126").Value = 10000
$ws_ARM.Range("I126").Value = 10000
$ws_ARM.Range("K126").Value = 30000
$ws_ARM.Range("M126").Value = -27530
$ws_ARM.Range("H136").Value = 3649.7778
$ws_ARM.Range("I136").Value = 3430.7144
$ws_ARM.Range("K136").Value = 10292.1432
$ws_ARM.Range("M136").Value = -7742.143199999999

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H132").Value = 115593.2
$ws_BSM.Range("J132").Value = 115593.2
$ws_BSM.Range("L132").Value = 115593.2
$ws_BSM.Range("N132").Value = -125713.2
$ws_BSM.Range("H135").Value = 0
$ws_BSM.Range("J135").Value = 0
$ws_BSM.Range("L135").Value = 0
$ws_BSM.Range("N135").ClearContents()

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H14").Value = 465.66666
$ws_CUL.Range("I14").Value = 465.66666
$ws_CUL.Range("K14").Value = 1396.99998
$ws_CUL.Range("M14").Value = -1223.99998
$ws_CUL.Range("H17").Value = 558.7273
$ws_CUL.Range("J17").Value = 580.4
$ws_CUL.Range("L17").Value = 1741.2
$ws_CUL.Range("N17").Value = -2079.2
$ws_CUL.Range("H18").Value = 2439.75
$ws_CUL.Range("I18").Value = 103.8
$ws_CUL.Range("K18").Value = 311.4
$ws_CUL.Range("M18").Value = -142.4
$ws_CUL.Range("H80").Value = 0
$ws_CUL.Range("J80").Value = 0
$ws_CUL.Range("L80").Value = 0
$ws_CUL.Range("N80").ClearContents()
$ws_CUL.Range("H83").Value = 0
$ws_CUL.Range("J83").Value = 0
$ws_CUL.Range("L83").Value = 0
$ws_CUL.Range("N83").ClearContents()
$ws_CUL.Range("H92").Value = 1249.75
$ws_CUL.Range("I92").Value = 0
$ws_CUL.Range("J92").Value = 1249.75
$ws_CUL.Range("K92").Value = 0
$ws_CUL.Range("L92").Value = 3749.25
$ws_CUL.Range("M92").ClearContents()
$ws_CUL.Range("N92").Value = -6245.25
$ws_CUL.Range("H95").Value = 9666.666999999999
$ws_CUL.Range("I95").Value = 8000
$ws_CUL.Range("J95").Value = 10000
$ws_CUL.Range("K95").Value = 24000
$ws_CUL.Range("L95").Value = 30000
$ws_CUL.Range("M95").Value = -21941
$ws_CUL.Range("N95").Value = -34118
$ws_CUL.Range("H107").Value = 799.93335
$ws_CUL.Range("J107").Value = 836.0278
$ws_CUL.Range("L107").Value = 2508.0834
$ws_CUL.Range("N107").Value = -6348.0834
$ws_CUL.Range("H113").Value = 1234.7222
$ws_CUL.Range("J113").Value = 1178.5186
$ws_CUL.Range("L113").Value = 3535.5558
$ws_CUL.Range("N113").Value = -7875.5558
$ws_CUL.Range("H122").Value = 2119.8
$ws_CUL.Range("I122").Value = 799.5
$ws_CUL.Range("K122").Value = 7195.5
$ws_CUL.Range("M122").Value = -4745.5
$ws_CUL.Range("H129").Value = 1109.1765
$ws_CUL.Range("I129").Value = 463.64285
$ws_CUL.Range("K129").Value = 1390.92855
$ws_CUL.Range("M129").Value = 3609.07145
$ws_CUL.Range("H131").Value = 6027
$ws_CUL.Range("J131").Value = 6082.5
$ws_CUL.Range("L131").Value = 18247.5
$ws_CUL.Range("N131").Value = -28327.5
$ws_CUL.Range("H132").Value = 2566.7334
$ws_CUL.Range("I132").Value = 2676.1667
$ws_CUL.Range("J132").Value = 2493.7778
$ws_CUL.Range("K132").Value = 24085.5003
$ws_CUL.Range("L132").Value = 22444.0002
$ws_CUL.Range("M132").Value = -21555.5003
$ws_CUL.Range("N132").Value = -27504.0002

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H70").Value = 3989.8572
$ws_GSM.Range("I70").Value = 2835.6667
$ws_GSM.Range("J70").Value = 4412.122
$ws_GSM.Range("K70").Value = 2835.6667
$ws_GSM.Range("L70").Value = 4412.122
$ws_GSM.Range("M70").Value = -2565.6667
$ws_GSM.Range("N70").Value = -4952.122
$ws_GSM.Range("H73").Value = 3989.8572
$ws_GSM.Range("I73").Value = 2835.6667
$ws_GSM.Range("J73").Value = 4412.122
$ws_GSM.Range("K73").Value = 2835.6667
$ws_GSM.Range("L73").Value = 4412.122
$ws_GSM.Range("M73").Value = -1899.6667
$ws_GSM.Range("N73").Value = -6284.122
$ws_GSM.Range("H97").Value = 982.8461
$ws_GSM.Range("I97").Value = 631.13635
$ws_GSM.Range("K97").Value = 631.13635
$ws_GSM.Range("M97").Value = -135.13635
$ws_GSM.Range("H102").Value = 2442.6875
$ws_GSM.Range("I102").Value = 2408.0667
$ws_GSM.Range("K102").Value = 2408.0667
$ws_GSM.Range("M102").Value = -786.0666999999999
$ws_GSM.Range("H122").Value = 4624.222
$ws_GSM.Range("I122").Value = 5014.75
$ws_GSM.Range("J122").Value = 1500
$ws_GSM.Range("K122").Value = 15044.25
$ws_GSM.Range("L122").Value = 4500
$ws_GSM.Range("M122").Value = -12594.25
$ws_GSM.Range("N122").Value = -9400
$ws_GSM.Range("H136").Value = 50625.344
$ws_GSM.Range("J136").Value = 50161.645
$ws_GSM.Range("L136").Value = 150484.935
$ws_GSM.Range("N136").Value = -155584.935

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H6").Value = 60069.668
$ws_LTW.Range("J6").Value = 60069.668
$ws_LTW.Range("L6").Value = 60069.668
$ws_LTW.Range("N6").Value = -60293.668
$ws_LTW.Range("H7").Value = 2387.6667
$ws_LTW.Range("I7").Value = 2265.2
$ws_LTW.Range("J7").Value = 3000
$ws_LTW.Range("K7").Value = 2265.2
$ws_LTW.Range("L7").Value = 3000
$ws_LTW.Range("M7").Value = -2153.2
$ws_LTW.Range("N7").Value = -3224
$ws_LTW.Range("H40").Value = 4488.7085
$ws_LTW.Range("I40").Value = 4644.2383
$ws_LTW.Range("J40").Value = 3400
$ws_LTW.Range("K40").Value = 4644.2383
$ws_LTW.Range("L40").Value = 3400
$ws_LTW.Range("M40").Value = -4508.2383
$ws_LTW.Range("N40").Value = -3672
$ws_LTW.Range("H55").Value = 971.5263
$ws_LTW.Range("I55").Value = 1057.091
$ws_LTW.Range("J55").Value = 853.875
$ws_LTW.Range("K55").Value = 1057.091
$ws_LTW.Range("L55").Value = 853.875
$ws_LTW.Range("M55").Value = -884.0909999999999
$ws_LTW.Range("N55").Value = -1199.875
$ws_LTW.Range("H100").Value = 1627.909
$ws_LTW.Range("I100").Value = 1627.909
$ws_LTW.Range("K100").Value = 1627.909
$ws_LTW.Range("M100").Value = -1086.909
$ws_LTW.Range("H122").Value = 4929180.5
$ws_LTW.Range("J122").Value = 3584.1667
$ws_LTW.Range("L122").Value = 10752.5001
$ws_LTW.Range("N122").Value = -15652.5001
$ws_LTW.Range("H126").Value = 2387.6667
$ws_LTW.Range("I126").Value = 2265.2
$ws_LTW.Range("J126").Value = 3000
$ws_LTW.Range("K126").Value = 6795.599999999999
$ws_LTW.Range("L126").Value = 9000
$ws_LTW.Range("M126").Value = -4325.599999999999
$ws_LTW.Range("N126").Value = -13940
$ws_LTW.Range("H128").Value = 99165
$ws_LTW.Range("J128").Value = 99165
$ws_LTW.Range("L128").Value = 99165
$ws_LTW.Range("N128").Value = -109125
$ws_LTW.Range("H132").Value = 3539.74
$ws_LTW.Range("I132").Value = 2657.25
$ws_LTW.Range("J132").Value = 3818.4211
$ws_LTW.Range("K132").Value = 7971.75
$ws_LTW.Range("L132").Value = 11455.2633
$ws_LTW.Range("M132").Value = -5441.75
$ws_LTW.Range("N132").Value = -16515.2633

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H96").Value = 2281.9167
$ws_WVR.Range("I96").Value = 2281.9167
$ws_WVR.Range("K96").Value = 2281.9167
$ws_WVR.Range("M96").Value = -908.9167000000002
$ws_WVR.Range("H132").Value = 3095.1843
$ws_WVR.Range("I132").Value = 2628
$ws_WVR.Range("K132").Value = 7884
$ws_WVR.Range("M132").Value = -5354
$ws_WVR.Range("H136").Value = 27947.05
$ws_WVR.Range("I136").Value = 2083.8965
$ws_WVR.Range("K136").Value = 6251.689499999999
$ws_WVR.Range("M136").Value = -3701.689499999999
